# Update packages for Document Understanding template.
# Remove the now-unused "STATE OF S TEXAS 3H1" value from D2 on both sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("D2")
    if ($cell.Value2 -eq "STATE OF S TEXAS 3H1") {
        $cell.ClearContents()
    }
}
